# Updated symbol list: refresh Price column (D) with latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    "D2" = "279.11"
    "D3" = "23.44"
    "D4" = "6.392"
    "D5" = "0.06284"
    "D6" = "3.632"
    "D7" = "6.620"
    "D8" = "1.394"
    "D9" = "0.8306"
    "D10" = "0.01389"
    "D11" = "0.1594"
    "D12" = "0.08405"
    "D13" = "0.03463"
    "D15" = "4.046"
    "D16" = "0.09294"
    "D17" = "0.001652"
    "D18" = "0.04743"
    "D19" = "0.006297"
    "D20" = "0.005945"
    "D21" = "0.001075"
    "D22" = "0.0001494"
    "D23" = "3.729"
    "D24" = "2.326"
    "D25" = "0.3334"
    "D26" = "0.1259"
    "D28" = "0.0002693"
    "D40" = "0.04766"
    "D41" = "0.007083"
    "D42" = "0.1175"
    "D43" = "0.003641"
    "D44" = "0.01233"
    "D45" = "0.00006069"
    "D46" = "0.0009863"
    "D48" = "0.7791"
    "D49" = "0.002468"
    "D51" = "0.01235"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $ws.Range($cellRef).Value = $priceUpdates[$cellRef]
}
